# Applies the "1st changes of mifos to finflux" edit:
#  - On the "Repayment schedule" sheet, a new (blank) column is inserted
#    before column N, shifting the old "Late" / "heading" (Disbursement) /
#    "Outstanding" columns (N, O, P) one place to the right (O, P, Q).
#  - The new column N takes on the width that column M (the column to its
#    left) already had.
#  - The "Repayment schedule" sheet becomes the active sheet/tab (it was
#    "Transactions" before), with the selection on that sheet set to S7.

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column); this shifts
# columns N:P -> O:Q, exactly as Excel's own "Insert Column" command does.
$schedule.Columns("N").Insert()

# The newly inserted column inherits the width of the column to its left.
$schedule.Columns("N").ColumnWidth = $schedule.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select cell S7 on it
# (previously "Transactions" was the active tab).
$schedule.Activate()
$schedule.Range("S7").Select() | Out-Null
